$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "48.436.11"
$ws.Range("E2").Value = "  +0.55%  "

# Row 3
$ws.Range("D3").Value = "2.502.11"
$ws.Range("E3").Value = "  -0.30%  "

# Row 4
$ws.Range("E4").Value = "  -0.17%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "318.96"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.45%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "106.23"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.79%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.521"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.60%  "

# Row 8
$ws.Range("E8").Value = "  -0.05%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.540"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.31%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.23"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.94%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.12"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.19%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0804"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.81%  "

# Row 13
$ws.Range("E13").Value = "  +1.09%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.10"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.46%  "

# Row 15
$ws.Range("D15").Value = "2.890.79"
$ws.Range("E15").Value = "  -0.40%  "

# Row 16
$ws.Range("D16").Value = "2.474.73"
$ws.Range("E16").Value = "  -1.48%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.841"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.87%  "

# Row 18
$ws.Range("D18").Value = "48.236.34"
$ws.Range("E18").Value = "  +0.46%  "

# Row 19
$ws.Range("B19").Value = "InternetComputer(DFINITY)"
$ws.Range("C19").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.80"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.14%  "

# Row 20
$ws.Range("B20").Value = "ImmutableX"
$ws.Range("C20").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "2.92"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +7.35%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.56"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.54%  "

# Row 22
$ws.Range("D22").Value = "0.0₃0934"
$ws.Range("E22").Value = "  -1.11%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "280.04"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.03%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "70.92"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.65%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.50"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.20%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "25.79"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.53%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.22"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -7.73%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.66"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.25%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.140"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.80%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "34.86"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.85%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "49.23"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.24%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "19.45"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.31%  "

# Row 34
$ws.Range("E34").Value = "  -0.29%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.27"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.47%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0772"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.35%  "

# Row 37
$ws.Range("E37").Value = "  -0.66%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.53"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.53%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.90"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.94%  "

# Row 40
$ws.Range("E40").Value = "  -0.90%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.21"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.01%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "119.90"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.64%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "21.81"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.77%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0300"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.77%  "

# Row 45
$ws.Range("D45").Value = "1.991.44"
$ws.Range("E45").Value = "  -1.97%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.17"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.56%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.97"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +5.66%  "

# Row 48
$ws.Range("E48").Value = "  +6.00%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.99"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.35%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "5.15"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.47%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "79.59"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.01%  "
